{"js": "// Update the multiplication-problem cells in the practice table.\n// Each \"AAA\u00d7B=\" expression in the document is unique, so a direct\n// search-and-replace per pair is unambiguous and safe.\nconst replacements = [\n  [\"668\u00d76=\", \"383\u00d74=\"],\n  [\"910\u00d72=\", \"560\u00d73=\"],\n  [\"227\u00d74=\", \"341\u00d75=\"],\n  [\"242\u00d78=\", \"396\u00d74=\"],\n  [\"232\u00d74=\", \"205\u00d73=\"],\n  [\"217\u00d74=\", \"801\u00d79=\"],\n  [\"993\u00d76=\", \"301\u00d75=\"],\n  [\"468\u00d75=\", \"208\u00d73=\"],\n  [\"751\u00d75=\", \"999\u00d79=\"],\n  [\"664\u00d72=\", \"651\u00d78=\"],\n  [\"495\u00d72=\", \"937\u00d77=\"],\n  [\"769\u00d78=\", \"113\u00d77=\"],\n  [\"675\u00d76=\", \"622\u00d73=\"],\n  [\"640\u00d75=\", \"336\u00d79=\"],\n  [\"994\u00d78=\", \"558\u00d79=\"],\n  [\"161\u00d74=\", \"601\u00d76=\"],\n  [\"669\u00d73=\", \"574\u00d79=\"],\n  [\"433\u00d79=\", \"356\u00d78=\"],\n  [\"224\u00d79=\", \"808\u00d72=\"],\n  [\"295\u00d79=\", \"803\u00d75=\"],\n  [\"708\u00d76=\", \"885\u00d72=\"],\n  [\"972\u00d79=\", \"856\u00d78=\"],\n  [\"207\u00d77=\", \"948\u00d78=\"],\n  [\"826\u00d75=\", \"164\u00d73=\"],\n  [\"636\u00d76=\", \"648\u00d75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the multiplication-problem cells in the practice table.\n# Each \"AAA\u00d7B=\" expression in the document is unique, so a direct\n# search-and-replace per pair is unambiguous and safe.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"668\u00d76=\", \"383\u00d74=\"),\n    @(\"910\u00d72=\", \"560\u00d73=\"),\n    @(\"227\u00d74=\", \"341\u00d75=\"),\n    @(\"242\u00d78=\", \"396\u00d74=\"),\n    @(\"232\u00d74=\", \"205\u00d73=\"),\n    @(\"217\u00d74=\", \"801\u00d79=\"),\n    @(\"993\u00d76=\", \"301\u00d75=\"),\n    @(\"468\u00d75=\", \"208\u00d73=\"),\n    @(\"751\u00d75=\", \"999\u00d79=\"),\n    @(\"664\u00d72=\", \"651\u00d78=\"),\n    @(\"495\u00d72=\", \"937\u00d77=\"),\n    @(\"769\u00d78=\", \"113\u00d77=\"),\n    @(\"675\u00d76=\", \"622\u00d73=\"),\n    @(\"640\u00d75=\", \"336\u00d79=\"),\n    @(\"994\u00d78=\", \"558\u00d79=\"),\n    @(\"161\u00d74=\", \"601\u00d76=\"),\n    @(\"669\u00d73=\", \"574\u00d79=\"),\n    @(\"433\u00d79=\", \"356\u00d78=\"),\n    @(\"224\u00d79=\", \"808\u00d72=\"),\n    @(\"295\u00d79=\", \"803\u00d75=\"),\n    @(\"708\u00d76=\", \"885\u00d72=\"),\n    @(\"972\u00d79=\", \"856\u00d78=\"),\n    @(\"207\u00d77=\", \"948\u00d78=\"),\n    @(\"826\u00d75=\", \"164\u00d73=\"),\n    @(\"636\u00d76=\", \"648\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $old\n    $range.Find.Replacement.Text = $new\n    $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
